$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old "Total Loan Amount / Booking Month / %As of Today X" helper table
# that lived in columns I:R (plus the merged captions in J3 and K5) is no
# longer needed - wipe its contents while leaving the cell formatting in
# place. (ClearContents doesn't reliably touch cells that are part of a
# merged range in this engine, so assign $null instead.)
$ws.Range("J3").Value = $null
$ws.Range("K5").Value = $null
$ws.Range("I6:R20").Value = $null

# Reset the view: select F4 (this also drops the old scrolled-down
# topLeftCell/selection that pointed at M106).
$ws.Range("F4").Select()
